$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("parameter_values")

# First enter the "squeeze_factor..." parameter name/value in row 84
$ws.Range("A84").Value = "squeeze_factor_threshold_delivery_attendance"
$ws.Range("B84").Value = 0.8

# Then insert a new row above it (pushing squeeze_factor down to row 85)
# and populate the new row 84 with the "dummy_prob_health_centre" parameter name/value
$ws.Rows.Item(84).Insert()
$ws.Range("A84").Value = "dummy_prob_health_centre"
$ws.Range("B84").Value = 0.7

# Finally, fill in the "DUMMY" source/comment column for both new rows
$ws.Range("C84").Value = "DUMMY"
$ws.Range("C85").Value = "DUMMY"

# Match the author's final view/selection state on the sheet
$ws.Activate() | Out-Null
$excel.ActiveWindow.ScrollRow = 73
$ws.Range("B87").Select() | Out-Null
